$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the record to remove by its unique MATRICULA value (column C),
# then delete that entire row, shifting the rows below it up - matching
# the data-cleanup performed by the authoring script.
$target = $ws.Cells.Find("GO338022654")
if ($target -ne $null) {
    $ws.Rows.Item($target.Row).Delete()
}
